# C5-PowerPoint.pptx edit
#
# 1) Slide 6 ("SOURCES OF FINANCE") table: switch the applied table style
#    from the deck's custom "Table_0" style to the built-in
#    "Medium Style 2 - Accent 1" style ({8B354DD5-45BE-4DEB-B63D-9CF03827D601}).
#
# 2) Re-colour the deck's theme (used by the slide master / all slides) from
#    the "Integral" palette over to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 -------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{8B354DD5-45BE-4DEB-B63D-9CF03827D601}", $true)

# --- 2) Theme colours: Integral -> Office ---------------------------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme

function Set-ThemeColor($index, $r, $g, $b) {
    $color = $themeColors.Colors($index)
    $color.RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeColor 1  0x00 0x00 0x00   # dk1
Set-ThemeColor 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor 12 0x95 0x4F 0x72   # folHlink
